$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: wipe all existing content so the old shared-strings table is
# fully dereferenced and gets rebuilt from scratch, in the order we write it.
$ws.Cells.Clear()

# Column letters for B..T (header row / value columns) and C..T (HKL / 1-columns)
$colsBtoT = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$colsCtoT = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# --- Step 2: Row 1 header numbers 0..18 in B1:T1
for ($i = 0; $i -lt $colsBtoT.Length; $i++) {
    $cell = $ws.Range($colsBtoT[$i] + "1")
    $cell.Value2 = $i
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- Step 3: Row 2 - A2 = 0 (header style), B2 = "HKL" (header style),
# C2:T2 = HKL index labels (normal style, string)
$a2 = $ws.Range("A2")
$a2.Value2 = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$b2 = $ws.Range("B2")
$b2.Value2 = "HKL"
$b2.Font.Bold = $true
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4160
$b2.Borders.LineStyle = 1

$hklLabels = @("[4, 0, 0]","[2, 1, 1]","[2, 2, 0]","[2, 0, 0]","[2, 2, 2]","[3, 1, 0]","[1, 1, 0]","[3, 2, 1]","1Pair-A","1Pair-B","2Pairs-A","2Pairs-B","3Pairs-A","3Pairs-B","3Pairs-C","4Pairs","5A4F","MaxUnique")
for ($i = 0; $i -lt $colsCtoT.Length; $i++) {
    $ws.Range($colsCtoT[$i] + "2").Value2 = $hklLabels[$i]
}

# --- Step 4: Rows 3..29 - A column index number (header style), B column
# label (header style), C:T = 1 (18 ones, normal style)
$rowLabels = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

$r = 3
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $aCell = $ws.Range("A" + $r)
    $aCell.Value2 = $i + 1
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $bCell = $ws.Range("B" + $r)
    $bCell.Value2 = $rowLabels[$i]

    foreach ($col in $colsCtoT) {
        $ws.Range($col + $r).Value2 = 1
    }

    $r++
}

$wb.Save()
